# Update scheduling figures for Oaks Crisis (jan_2021 column K and the
# SFY 2021 Total column Q) on the crisis_src worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("crisis_src")

# Map of cell address -> new value
$updates = @{
    "K3"  = 289
    "Q3"  = 1871
    "Q4"  = 256
    "K5"  = 346
    "Q5"  = 2172
    "K8"  = 41
    "Q8"  = 219
    "K9"  = 3
    "Q9"  = 16
    "K10" = 12
    "Q10" = 40
    "K13" = 27
    "Q13" = 83
    "K14" = 83
    "Q14" = 360
    "K18" = 27
    "Q18" = 106
    "K20" = 10
    "Q20" = 48
    "K21" = 56
    "Q21" = 236
    "K22" = 4
    "Q22" = 10
    "K23" = 2
    "Q23" = 9
    "K24" = 6
    "Q24" = 12
    "K26" = 5
    "Q26" = 12
    "K28" = 26
    "Q28" = 124
    "K29" = 13
    "Q29" = 65
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
